# =============================================================================
# Add 2022-Q4 data (new quarterly snapshot inserted right after "总计").
# =============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for 2022-Q4 at row 2, pushing
#    the existing quarters down by one row (column A is a running 0-based
#    index, so it is renumbered after the shift).
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $srcRange = $totalWs.Range("A" + $r + ":D" + $r)
    $dstRange = $totalWs.Range("A" + $dst + ":D" + $dst)
    $srcRange.Copy($dstRange)
}

$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 10
$totalWs.Range("D2").Value = 4.73

for ($r = 2; $r -le 9; $r++) {
    $totalWs.Cells.Item($r,1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Insert a new "2022-Q4" worksheet right after "总计", holding the fund
#    holdings data, formatted like the other quarterly sheets.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $totalWs)
$ws.Name = "2022-Q4"

# Copy header formatting (same Chinese column titles/style) from "2022-Q3".
$tmpl = $wb.Worksheets.Item("2022-Q3")
$tmpl.Range("B1:H1").Copy($ws.Range("B1:H1"))

# Copy column-A numeric style (bold/border/center like the other quarter
# sheets use for their row-index column).
$tmpl.Range("A2:A11").Copy($ws.Range("A2:A11"))

# Fill in fund-holdings data rows 2-11 (columns B and D-G keep their literal
# text formatting -- leading zeros / fixed decimals -- instead of becoming
# numbers, matching how the other quarterly sheets store this data).
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'003567"
$ws.Cells.Item(2,2).Style = "Normal"
$ws.Cells.Item(2,3).Value = "华夏行业景气混合"
$ws.Cells.Item(2,4).Value = "'109.60"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'93.65"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(2,6).Value = "'2.85"
$ws.Cells.Item(2,6).Style = "Normal"
$ws.Cells.Item(2,7).Value = "'3.1236"
$ws.Cells.Item(2,7).Style = "Normal"
$ws.Cells.Item(2,8).Value = 3
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'001349"
$ws.Cells.Item(3,2).Style = "Normal"
$ws.Cells.Item(3,3).Value = "富国改革动力混合"
$ws.Cells.Item(3,4).Value = "'15.49"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'89.34"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(3,6).Value = "'3.43"
$ws.Cells.Item(3,6).Style = "Normal"
$ws.Cells.Item(3,7).Value = "'0.5313"
$ws.Cells.Item(3,7).Style = "Normal"
$ws.Cells.Item(3,8).Value = 5
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'003378"
$ws.Cells.Item(4,2).Style = "Normal"
$ws.Cells.Item(4,3).Value = "泰康策略优选灵活配置混合"
$ws.Cells.Item(4,4).Value = "'13.78"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'83.09"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(4,6).Value = "'2.84"
$ws.Cells.Item(4,6).Style = "Normal"
$ws.Cells.Item(4,7).Value = "'0.3914"
$ws.Cells.Item(4,7).Style = "Normal"
$ws.Cells.Item(4,8).Value = 9
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'952004"
$ws.Cells.Item(5,2).Style = "Normal"
$ws.Cells.Item(5,3).Value = "国泰君安君得明混合"
$ws.Cells.Item(5,4).Value = "'18.25"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'83.43"
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(5,6).Value = "'1.41"
$ws.Cells.Item(5,6).Style = "Normal"
$ws.Cells.Item(5,7).Value = "'0.2573"
$ws.Cells.Item(5,7).Style = "Normal"
$ws.Cells.Item(5,8).Value = 10
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'011769"
$ws.Cells.Item(6,2).Style = "Normal"
$ws.Cells.Item(6,3).Value = "富国精诚回报12个月持有期混合A"
$ws.Cells.Item(6,4).Value = "'21.97"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'24.85"
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(6,6).Value = "'1.08"
$ws.Cells.Item(6,6).Style = "Normal"
$ws.Cells.Item(6,7).Value = "'0.2373"
$ws.Cells.Item(6,7).Style = "Normal"
$ws.Cells.Item(6,8).Value = 6
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'952035"
$ws.Cells.Item(7,2).Style = "Normal"
$ws.Cells.Item(7,3).Value = "国泰君安君得诚混合"
$ws.Cells.Item(7,4).Value = "'2.21"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'85.83"
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(7,6).Value = "'3.25"
$ws.Cells.Item(7,6).Style = "Normal"
$ws.Cells.Item(7,7).Value = "'0.0718"
$ws.Cells.Item(7,7).Style = "Normal"
$ws.Cells.Item(7,8).Value = 9
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'010495"
$ws.Cells.Item(8,2).Style = "Normal"
$ws.Cells.Item(8,3).Value = "创金合信创新驱动股票A"
$ws.Cells.Item(8,4).Value = "'0.93"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'90.74"
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(8,6).Value = "'7.60"
$ws.Cells.Item(8,6).Style = "Normal"
$ws.Cells.Item(8,7).Value = "'0.0707"
$ws.Cells.Item(8,7).Style = "Normal"
$ws.Cells.Item(8,8).Value = 1
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'010496"
$ws.Cells.Item(9,2).Style = "Normal"
$ws.Cells.Item(9,3).Value = "创金合信创新驱动股票C"
$ws.Cells.Item(9,4).Value = "'0.32"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'90.74"
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(9,6).Value = "'7.60"
$ws.Cells.Item(9,6).Style = "Normal"
$ws.Cells.Item(9,7).Value = "'0.0243"
$ws.Cells.Item(9,7).Style = "Normal"
$ws.Cells.Item(9,8).Value = 1
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'011770"
$ws.Cells.Item(10,2).Style = "Normal"
$ws.Cells.Item(10,3).Value = "富国精诚回报12个月持有期混合C"
$ws.Cells.Item(10,4).Value = "'1.39"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'24.85"
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(10,6).Value = "'1.08"
$ws.Cells.Item(10,6).Style = "Normal"
$ws.Cells.Item(10,7).Value = "'0.0150"
$ws.Cells.Item(10,7).Style = "Normal"
$ws.Cells.Item(10,8).Value = 6
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'000788"
$ws.Cells.Item(11,2).Style = "Normal"
$ws.Cells.Item(11,3).Value = "前海开源中国成长灵活配置混合"
$ws.Cells.Item(11,4).Value = "'0.43"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'85.88"
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(11,6).Value = "'1.60"
$ws.Cells.Item(11,6).Style = "Normal"
$ws.Cells.Item(11,7).Value = "'0.0069"
$ws.Cells.Item(11,7).Style = "Normal"
$ws.Cells.Item(11,8).Value = 9

# Match the page margins used throughout the rest of the workbook.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
